$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "62.003.29"
$ws.Range("E2").Value = "  +1.78%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.473.80"
$ws.Range("E3").Value = "  +3.06%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5: BNB
$ws.Range("D5").Value = "'583.37"
$ws.Range("E5").Value = "  +1.96%  "

# Row 6: Solana
$ws.Range("D6").Value = "'148.54"
$ws.Range("E6").Value = "  +9.34%  "

# Row 7: LidoStakedEther
$ws.Range("D7").Value = "3.473.19"
$ws.Range("E7").Value = "  +3.10%  "

# Row 8: USDC
$ws.Range("E8").Value = "  +0.04%  "

# Row 9: XRP
$ws.Range("D9").Value = "'0.477"
$ws.Range("E9").Value = "  +1.76%  "

# Row 10: Toncoin
$ws.Range("D10").Value = "'7.73"

# Row 11: Dogecoin
$ws.Range("D11").Value = "'0.127"
$ws.Range("E11").Value = "  +4.17%  "

# Row 12: Cardano
$ws.Range("D12").Value = "'0.395"
$ws.Range("E12").Value = "  +3.09%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.067.87"
$ws.Range("E13").Value = "  +3.04%  "

# Row 14: TRON->Avalanche
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'27.95"
$ws.Range("E14").Value = "  +10.08%  "

# Row 15: Avalanche->TRON
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "'0.122"
$ws.Range("E15").Value = "  -0.67%  "

# Row 16: ShibaInu
$ws.Range("D16").Value = "'0.0000176"
$ws.Range("E16").Value = "  +2.48%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "3.482.56"
$ws.Range("E17").Value = "  +3.11%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "62.152.88"
$ws.Range("E18").Value = "  +1.64%  "

# Row 19: Polkadot
$ws.Range("D19").Value = "'6.29"
$ws.Range("E19").Value = "  +9.51%  "

# Row 20: Chainlink
$ws.Range("D20").Value = "'14.46"
$ws.Range("E20").Value = "  +4.92%  "

# Row 21: Uniswap
$ws.Range("D21").Value = "'9.59"
$ws.Range("E21").Value = "  +3.28%  "

# Row 22: BitcoinCash
$ws.Range("D22").Value = "'390.87"
$ws.Range("E22").Value = "  +4.63%  "

# Row 23: Polygon
$ws.Range("D23").Value = "'0.568"
$ws.Range("E23").Value = "  +3.74%  "

# Row 24: Litecoin
$ws.Range("D24").Value = "'73.64"
$ws.Range("E24").Value = "  +3.99%  "

# Row 25: Dai->LEO
$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").Value = "'5.78"
$ws.Range("E25").Value = "  +0.30%  "

# Row 26: LEO->Dai
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.07%  "

# Row 27: PEPE
$ws.Range("D27").Value = "'0.0000125"
$ws.Range("E27").Value = "  -0.23%  "

# Row 28: WrappedeETH
$ws.Range("D28").Value = "3.600.24"
$ws.Range("E28").Value = "  +2.40%  "

# Row 29: Kaspa
$ws.Range("D29").Value = "'0.181"
$ws.Range("E29").Value = "  -0.05%  "

# Row 30: RenderToken
$ws.Range("D30").Value = "'7.82"
$ws.Range("E30").Value = "  +5.76%  "

# Row 31: Binance-PegBSC-USD
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.41%  "

# Row 32: Fetch.AI
$ws.Range("E32").Value = "  -8.42%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").Value = "'8.28"
$ws.Range("E33").Value = "  +2.90%  "

# Row 34: PancakeSwap
$ws.Range("D34").Value = "'2.20"
$ws.Range("E34").Value = "  +3.01%  "

# Row 35: USDe
$ws.Range("E35").Value = "  +0.02%  "

# Row 36: EthereumClassic
$ws.Range("D36").Value = "'24.35"
$ws.Range("E36").Value = "  +3.77%  "

# Row 37: RenzoRestakedETH
$ws.Range("D37").Value = "3.503.54"
$ws.Range("E37").Value = "  +3.18%  "

# Row 38: Aptos
$ws.Range("D38").Value = "'7.06"
$ws.Range("E38").Value = "  +4.01%  "

# Row 39: NEARProtocol
$ws.Range("D39").Value = "'5.22"
$ws.Range("E39").Value = "  +0.67%  "

# Row 40: ImmutableX
$ws.Range("D40").Value = "'1.58"
$ws.Range("E40").Value = "  +2.79%  "

# Row 41: Monero
$ws.Range("D41").Value = "'166.55"
$ws.Range("E41").Value = "  +1.18%  "

# Row 42: Hedera
$ws.Range("D42").Value = "'0.0790"
$ws.Range("E42").Value = "  +4.90%  "

# Row 43: EnergySwap
$ws.Range("D43").Value = "'27.62"
$ws.Range("E43").Value = "  +11.12%  "

# Row 44: Mantle
$ws.Range("D44").Value = "'0.811"
$ws.Range("E44").Value = "  +5.03%  "

# Row 45: FirstDigitalUSD->OKB
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'42.63"
$ws.Range("E45").Value = "  +2.10%  "

# Row 46: Filecoin
$ws.Range("D46").Value = "'4.52"
$ws.Range("E46").Value = "  +4.11%  "

# Row 47: OKB->Stacks
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'1.74"
$ws.Range("E47").Value = "  +2.84%  "

# Row 48: Stacks->FirstDigitalUSD
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  -0.01%  "

# Row 49: ONDO
$ws.Range("D49").Value = "'1.18"
$ws.Range("E49").Value = "  -0.49%  "

# Row 50: Maker
$ws.Range("D50").Value = "2.579.73"
$ws.Range("E50").Value = "  +4.71%  "

# Row 51: Cosmos
$ws.Range("D51").Value = "'6.95"
$ws.Range("E51").Value = "  +2.68%  "
